# correção das notas do fórum para matc65 em 2021.2
# Zera as colunas B (views diárias) até J (nota_view) para as linhas de
# alunos que tinham alguma visualização registrada (total_views > 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A área de dados vai da linha 2 até a linha 50 (cabeçalho na linha 1),
# colunas B (2022-09-18) até J (nota_view).
$ws.Range("B2:J50").Value = 0
